# Apply cryptocurrency price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.314.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "'3.505.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'584.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'134.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("D7").Value = "'3.506.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'7.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("D13").Value = "'4.101.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "'3.505.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'26.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.90%  "
$ws.Range("D18").Value = "'64.300.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "'9.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "'383.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "'3.641.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'74.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").Value = "'3.523.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").Value = "'164.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("D42").Value = "'0.0783"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "'26.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'41.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "'1.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "'2.480.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +0.82%  "
